$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Unity"
$ws.Range("C9").Value = "Unity"
$ws.Range("C12").Value = "Unity"
$ws.Range("C15").Value = "Unity"
$ws.Range("C22").Value = "Unity"
$ws.Range("E22").Value = "Done"

$ws.Range("D6").Select()
